# Update "江西-漫展信息.xlsx": refresh a handful of "want-to-go" counters
# and insert one brand-new exhibition row ("赣州·明日方舟only叙拉古夜宴3.0...")
# into both the "展览" sheet and the combined "全部类型" sheet.

$wb = $excel.ActiveWorkbook

function Update-Counts {
    param($ws, [hashtable]$updates)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item([int]$row, 6).Value = $updates[$row]
    }
}

function Set-TextValue {
    # Assign $value to the cell as literal text, even if it looks like a
    # date/number to Excel's input parser, and do so without leaving any
    # lingering NumberFormat/quotePrefix style on the cell.
    param($ws, [int]$row, [int]$col, $value, [int]$styleSourceRow, [int]$styleSourceCol)

    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $value
    $ws.Cells.Item($styleSourceRow, $styleSourceCol).Copy()
    $c.PasteSpecial(-4122)
}

function Insert-NewRow {
    param($ws, [int]$newRow, [int]$aboveRow, [int]$aNumber, [int]$lastRow)

    $ws.Rows.Item($newRow).Insert()

    # Column A: number + the same style (border/bold/center) as every other
    # data row, copied from the row immediately above the inserted one.
    $ws.Cells.Item($aboveRow, 1).Copy()
    $ws.Cells.Item($newRow, 1).PasteSpecial(-4122)
    $ws.Cells.Item($newRow, 1).Value = $aNumber

    Set-TextValue $ws $newRow 2 "2024-07-28" $aboveRow 2
    Set-TextValue $ws $newRow 3 "赣州·明日方舟only叙拉古夜宴3.0暨同好交流茶话会" $aboveRow 3
    Set-TextValue $ws $newRow 4 "兴国路恒大帝景西门 江西长庚控股有限公司" $aboveRow 4
    Set-TextValue $ws $newRow 5 "2024.07.28 11:00-07.28 17:00" $aboveRow 5

    $ws.Cells.Item($newRow, 6).Value = 0

    Set-TextValue $ws $newRow 7 "不可售" $aboveRow 7
    Set-TextValue $ws $newRow 8 "https://show.bilibili.com/platform/detail.html?id=85688" $aboveRow 8
    Set-TextValue $ws $newRow 9 "//i1.hdslb.com/bfs/openplatform/202405/5AFwM8QV1715765287721.png" $aboveRow 9

    # Every row the insert pushed down keeps its own "#" serial (column A)
    # one higher than it used to be reported as (serial = row - 1).
    for ($r = $newRow + 1; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 1
    }
}

# ---------------------------------------------------------------------------
# Sheet "展览" (exhibitions)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

Update-Counts $ws1 @{
    3  = 179   # was 178
    4  = 92    # was 90
    7  = 1732  # was 1728
    8  = 38    # was 37
    11 = 1839  # was 1832
    13 = 118   # was 116
    15 = 273   # was 272
    17 = 4     # was 2
    22 = 793   # was 788
}

Insert-NewRow $ws1 23 22 22 27

# The "南昌·幻梦境国际动漫游戏嘉年华1th" row (old row 23) slides down to row
# 24 as part of the insert above, and its "want to go" counter also ticked
# up from 312 to 314 in this same update.
$ws1.Cells.Item(24, 6).Value = 314

# ---------------------------------------------------------------------------
# Sheet "全部类型" (all event types)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

Update-Counts $ws4 @{
    3  = 179   # was 178
    4  = 92    # was 90
    7  = 1732  # was 1728
    9  = 38    # was 37
    12 = 1839  # was 1832
    14 = 118   # was 116
    16 = 273   # was 272
    18 = 4     # was 2
    23 = 793   # was 788
}

Insert-NewRow $ws4 24 23 23 28

# Same counter bump for "南昌·幻梦境国际动漫游戏嘉年华1th", which lives one
# row further down on this combined sheet (old row 24 -> new row 25).
$ws4.Cells.Item(25, 6).Value = 314
